# Update view/hit counts in column F for several rows across three sheets,
# reflecting a refreshed data export (gh-pages output at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 5329   # 南宁·AB动漫游戏嘉年华
$ws1.Range("F4").Value  = 11061  # 南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）
$ws1.Range("F5").Value  = 274    # 南宁·火影忍者only
$ws1.Range("F6").Value  = 581    # 南宁·蔚蓝档案only
$ws1.Range("F7").Value  = 156    # 南宁·国乙only
$ws1.Range("F8").Value  = 219    # 南宁·熊喵M动漫嘉年华【免费】
$ws1.Range("F9").Value  = 941    # 南宁·第二届北极光动漫展
$ws1.Range("F10").Value = 92     # 南宁·万圣漫控嘉年华10

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 22      # 南宁·跨越二次元ACG神级动漫世界巡回演唱会——

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 5329   # 南宁·AB动漫游戏嘉年华
$ws4.Range("F6").Value  = 22     # 南宁·跨越二次元ACG神级动漫世界巡回演唱会——
$ws4.Range("F7").Value  = 11061  # 南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）
$ws4.Range("F8").Value  = 274    # 南宁·火影忍者only
$ws4.Range("F9").Value  = 581    # 南宁·蔚蓝档案only
$ws4.Range("F10").Value = 156    # 南宁·国乙only
$ws4.Range("F13").Value = 219    # 南宁·熊喵M动漫嘉年华【免费】
$ws4.Range("F14").Value = 941    # 南宁·第二届北极光动漫展
$ws4.Range("F16").Value = 92     # 南宁·万圣漫控嘉年华10
